$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting of the last populated row (row 9) down across all new rows first,
# so the new rows inherit the same borders/alignment/number-format as the existing table
# without Excel fabricating brand-new style entries.
$ws.Range("A9:H9").Copy()
$ws.Range("A10:H32").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A10").Value = "CRS_REV_009"
$ws.Range("B10").Value = 44892
$ws.Range("C10").Value = "Esraa Abdelnaby"
$ws.Range("D10").Value = "V1.2"
$ws.Range("E10").Value = "page 4"
$ws.Range("F10").Value = "edit the context diagram"
$ws.Range("H10").Value = "Open"

$ws.Range("A11").Value = "CRS_REV_010"
$ws.Range("B11").Value = 44892
$ws.Range("C11").Value = "Esraa Abdelnaby"
$ws.Range("D11").Value = "V1.2"
$ws.Range("E11").Value = "page 5"
$ws.Range("F11").Value = "add requirments based on the change request of the CR"
$ws.Range("H11").Value = "Open"

$ws.Range("A12").Value = "CRS_REV_011"
$ws.Range("B12").Value = 44892
$ws.Range("C12").Value = "Esraa Abdelnaby"
$ws.Range("D12").Value = "V1.2"
$ws.Range("E12").Value = "page 5"
$ws.Range("F12").Value = "Allow sign up with `"Gmail`" & `"Fcebook`""
$ws.Range("H12").Value = "Open"

$ws.Range("A13").Value = "CRS_REV_012"
$ws.Range("B13").Value = 44892
$ws.Range("C13").Value = "Esraa Abdelnaby"
$ws.Range("D13").Value = "V1.2"
$ws.Range("E13").Value = "page 5"
$ws.Range("F13").Value = "add reset password option "
$ws.Range("H13").Value = "Open"

$ws.Range("A14").Value = "CRS_REV_013"
$ws.Range("B14").Value = 44892
$ws.Range("C14").Value = "Esraa Abdelnaby"
$ws.Range("D14").Value = "V1.2"
$ws.Range("E14").Value = "page 5"
$ws.Range("F14").Value = "Verify the accuracy of the information while logging in"
$ws.Range("H14").Value = "Open"

$ws.Range("A15").Value = "CRS_REV_014"
$ws.Range("B15").Value = 44892
$ws.Range("C15").Value = "Esraa Abdelnaby"
$ws.Range("D15").Value = "V1.2"
$ws.Range("E15").Value = "page 5"
$ws.Range("F15").Value = "Prevent the duplication of email during registration"
$ws.Range("H15").Value = "Open"

$ws.Range("A16").Value = "CRS_REV_015"
$ws.Range("B16").Value = 44892
$ws.Range("C16").Value = "Esraa Abdelnaby"
$ws.Range("D16").Value = "V1.2"
$ws.Range("E16").Value = "page 5"
$ws.Range("F16").Value = "app shall let users choose the plan "
$ws.Range("H16").Value = "Open"

$ws.Range("A17").Value = "CRS_REV_016"
$ws.Range("B17").Value = 44892
$ws.Range("C17").Value = "Esraa Abdelnaby"
$ws.Range("D17").Value = "V1.2"
$ws.Range("E17").Value = "page 6"
$ws.Range("F17").Value = "The ride starts after the current location is determined"
$ws.Range("H17").Value = "Open"

$ws.Range("A18").Value = "CRS_REV_017"
$ws.Range("B18").Value = 44892
$ws.Range("C18").Value = "Esraa Abdelnaby"
$ws.Range("D18").Value = "V1.2"
$ws.Range("E18").Value = "page 6"
$ws.Range("F18").Value = "The app shall automatically detect the current location if it has access to the location"
$ws.Range("H18").Value = "Open"

$ws.Range("A19").Value = "CRS_REV_018"
$ws.Range("B19").Value = 44892
$ws.Range("C19").Value = "Esraa Abdelnaby"
$ws.Range("D19").Value = "V1.2"
$ws.Range("E19").Value = "page 6"
$ws.Range("F19").Value = "app shall allow reporting of any issue "
$ws.Range("H19").Value = "Open"

$ws.Range("A20").Value = "CRS_REV_019"
$ws.Range("B20").Value = 44892
$ws.Range("C20").Value = "Esraa Abdelnaby"
$ws.Range("D20").Value = "V1.2"
$ws.Range("E20").Value = "page 7"
$ws.Range("F20").Value = "users of basic plan can upgrade to premium plan "
$ws.Range("H20").Value = "Open"

$ws.Range("A21").Value = "CRS_REV_020"
$ws.Range("B21").Value = 44892
$ws.Range("C21").Value = "Esraa Abdelnaby"
$ws.Range("D21").Value = "V1.2"
$ws.Range("E21").Value = "page 7"
$ws.Range("F21").Value = "on basic plan , the icons of all saved bumbs appear on the gps screen "
$ws.Range("H21").Value = "Open"

$ws.Range("A22").Value = "CRS_REV_021"
$ws.Range("B22").Value = 44892
$ws.Range("C22").Value = "Esraa Abdelnaby"
$ws.Range("D22").Value = "V1.2"
$ws.Range("E22").Value = "page 7"
$ws.Range("F22").Value = "voice notify users with the recorded speed bumbs "
$ws.Range("H22").Value = "Open"

$ws.Range("A23").Value = "CRS_REV_022"
$ws.Range("B23").Value = 44892
$ws.Range("C23").Value = "Esraa Abdelnaby"
$ws.Range("D23").Value = "V1.2"
$ws.Range("E23").Value = "page 7"
$ws.Range("F23").Value = "at choosing premium plan , app shall redirect the user to the website"
$ws.Range("H23").Value = "Open"

$ws.Range("A24").Value = "CRS_REV_023"
$ws.Range("B24").Value = 44892
$ws.Range("C24").Value = "Esraa Abdelnaby"
$ws.Range("D24").Value = "V1.2"
$ws.Range("E24").Value = "page 8"
$ws.Range("F24").Value = "on premium plan ,any upcoming bump shall appear on the screen"
$ws.Range("H24").Value = "Open"

$ws.Range("A25").Value = "CRS_REV_024"
$ws.Range("B25").Value = 44892
$ws.Range("C25").Value = "Esraa Abdelnaby"
$ws.Range("D25").Value = "V1.2"
$ws.Range("E25").Value = "page 8"
$ws.Range("F25").Value = "The app notifies the user of the time left until the next bump"
$ws.Range("H25").Value = "Open"

$ws.Range("A26").Value = "CRS_REV_025"
$ws.Range("B26").Value = 44893
$ws.Range("C26").Value = "Esraa Abdelnaby"
$ws.Range("D26").Value = "V1.2"
$ws.Range("E26").Value = "page 8"
$ws.Range("F26").Value = "by the end of the ride ,the app shall save the detected bumps"
$ws.Range("H26").Value = "Open"

$ws.Range("A27").Value = "CRS_REV_026"
$ws.Range("B27").Value = 44894
$ws.Range("C27").Value = "Esraa Abdelnaby"
$ws.Range("D27").Value = "V1.2"
$ws.Range("E27").Value = "page 8"
$ws.Range("F27").Value = "The user has to give the permission to connect to the detecting device via bluetooth"
$ws.Range("H27").Value = "Open"

$ws.Range("A28").Value = "CRS_REV_027"
$ws.Range("H28").Value = "Open"

$ws.Range("A29").Value = "CRS_REV_028"
$ws.Range("H29").Value = "Open"

$ws.Range("A30").Value = "CRS_REV_029"
$ws.Range("H30").Value = "Open"

$ws.Range("A31").Value = "CRS_REV_030"
$ws.Range("H31").Value = "Open"

$ws.Range("A32").Value = "CRS_REV_031"
$ws.Range("H32").Value = "Open"

$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("E13").Select()
